$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: header block (bold), copy of row 8's header ---
$ws.Range("A26").Value = "Milestone3"
$ws.Range("C26").Value = "gamma = 0.5"
$ws.Range("E26").Value = "alpha = 0.4"
$ws.Range("G26").Value = "144 states"
$ws.Range("H26").Value = "   13824 actions"
$ws.Range("A26:I26").Font.Bold = $true

# --- Row 27: action range labels ---
$ws.Range("A27").Value = "0-9000 actions 60% random"
$ws.Range("D27").Value = "9000-18000 actions 40% random"
$ws.Range("H27").Value = "18000-27000 actions 20% random"

# --- Row 28: fraction labels ---
$ws.Range("B28").Value = "1811/9000"
$ws.Range("D28").Value = "3585/9000"
$ws.Range("H28").Value = "4343/9000"

# --- Row 29: Training percentages ---
$ws.Range("A29").Value = "Training"
$ws.Range("B29").Value = 0.201
$ws.Range("D29").Value = 0.398
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = 0.483
$ws.Range("B29").NumberFormat = "0.00%"
$ws.Range("D29").NumberFormat = "0.00%"
$ws.Range("G29").NumberFormat = "0.00%"
$ws.Range("H29").NumberFormat = "0.00%"

# --- Row 30: Evaluation ---
$ws.Range("A30").Value = "Evaluation (0%)"
$ws.Range("C30").NumberFormat = "0.00%"

$ws.Range("E35").Select() | Out-Null
